# Auto-generated edit script: updates currency/price-related cells
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets to reflect refreshed
# market-board pricing data pulled by the scheduled runner.

$wb = $excel.ActiveWorkbook

# ================= Sheet: ALC =================
$ws = $wb.Worksheets.Item("ALC")

# Row 40
$ws.Range("H40").Value = 2437.625
$ws.Range("I40").Value = 2437.625
$ws.Range("K40").Value = 2437.625
$ws.Range("M40").Value = -2262.625
# Row 58
$ws.Range("H58").Value = 1301.52
$ws.Range("J58").Value = 1823.8667
$ws.Range("L58").Value = 5471.6001
$ws.Range("N58").Value = -5771.6001
# Row 62
$ws.Range("H62").Value = 55558556
$ws.Range("I62").Value = 111111110
$ws.Range("K62").Value = 111111110
$ws.Range("M62").Value = -111110486
# Row 65
$ws.Range("H65").Value = 55558556
$ws.Range("I65").Value = 111111110
$ws.Range("K65").Value = 555555550
$ws.Range("M65").Value = -555552430
# Row 70
$ws.Range("H70").Value = 986.1429000000001
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 986.1429000000001
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 2958.4287
$ws.Range("M70").ClearContents()
$ws.Range("N70").Value = -3498.4287
# Row 73
$ws.Range("H73").Value = 986.1429000000001
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 986.1429000000001
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 2958.4287
$ws.Range("M73").ClearContents()
$ws.Range("N73").Value = -4830.4287
# Row 76
$ws.Range("H76").Value = 5266.533
$ws.Range("J76").Value = 5533.1113
$ws.Range("L76").Value = 5533.1113
$ws.Range("N76").Value = -6163.1113
# Row 79
$ws.Range("H79").Value = 5266.533
$ws.Range("J79").Value = 5533.1113
$ws.Range("L79").Value = 5533.1113
$ws.Range("N79").Value = -7717.1113
# Row 86
$ws.Range("H86").Value = 3388.889
$ws.Range("I86").Value = 3660
$ws.Range("K86").Value = 3660
$ws.Range("M86").Value = -2537
# Row 89
$ws.Range("H89").Value = 3388.889
$ws.Range("I89").Value = 3660
$ws.Range("K89").Value = 18300
$ws.Range("M89").Value = -12684
# Row 106
$ws.Range("H106").Value = 9212.267
$ws.Range("I106").Value = 11471.272
$ws.Range("K106").Value = 11471.272
$ws.Range("M106").Value = -10840.272
# Row 112
$ws.Range("H112").Value = 2161.45
$ws.Range("I112").Value = 833
$ws.Range("J112").Value = 2395.8823
$ws.Range("K112").Value = 2499
$ws.Range("L112").Value = 7187.646900000001
$ws.Range("M112").Value = -1391
$ws.Range("N112").Value = -9403.6469
# Row 129
$ws.Range("H129").Value = 735.4706
$ws.Range("J129").Value = 925.9
$ws.Range("L129").Value = 2777.7
$ws.Range("N129").Value = -12777.7
# Row 137
$ws.Range("H137").Value = 1455.2941
$ws.Range("I137").Value = 1160.2727
$ws.Range("K137").Value = 3480.8181
$ws.Range("M137").Value = -930.8181
# Row 138
$ws.Range("H138").Value = 2361.4473
$ws.Range("I138").Value = 3875
$ws.Range("J138").Value = 2254.8591
$ws.Range("K138").Value = 11625
$ws.Range("L138").Value = 6764.577300000001
$ws.Range("M138").Value = -6485
$ws.Range("N138").Value = -17044.5773

# ================= Sheet: ARM =================
$ws = $wb.Worksheets.Item("ARM")

# Row 32
$ws.Range("H32").Value = 6260.7417
$ws.Range("I32").Value = 6322.7954
$ws.Range("K32").Value = 6322.7954
$ws.Range("M32").Value = -6035.7954
# Row 97
$ws.Range("H97").Value = 634
$ws.Range("I97").Value = 634
$ws.Range("K97").Value = 634
$ws.Range("M97").Value = -138
# Row 132
$ws.Range("H132").Value = 2533.9814
$ws.Range("I132").Value = 1820.9445
$ws.Range("J132").Value = 3960.0557
$ws.Range("K132").Value = 5462.833500000001
$ws.Range("L132").Value = 11880.1671
$ws.Range("M132").Value = -2932.833500000001
$ws.Range("N132").Value = -16940.1671

# ================= Sheet: BSM =================
$ws = $wb.Worksheets.Item("BSM")

# Row 20
$ws.Range("H20").Value = 1753.65
$ws.Range("I20").Value = 1859.1333
$ws.Range("K20").Value = 1859.1333
$ws.Range("M20").Value = -1612.1333
# Row 99
$ws.Range("H99").Value = 38462556
$ws.Range("I99").Value = 50000924
$ws.Range("J99").Value = 1325
$ws.Range("K99").Value = 50000924
$ws.Range("L99").Value = 1325
$ws.Range("M99").Value = -49999426
$ws.Range("N99").Value = -4321
# Row 134
$ws.Range("H134").Value = 928.6667
$ws.Range("I134").Value = 928.6667
$ws.Range("K134").Value = 2786.0001
$ws.Range("M134").Value = -251.0001000000002

# ================= Sheet: CRP =================
$ws = $wb.Worksheets.Item("CRP")

# Row 22
$ws.Range("H22").Value = 77967
$ws.Range("I22").Value = 141.83333
$ws.Range("J22").Value = 233617.33
$ws.Range("K22").Value = 141.83333
$ws.Range("L22").Value = 233617.33
$ws.Range("M22").Value = 208.16667
$ws.Range("N22").Value = -234317.33
# Row 31
$ws.Range("H31").Value = 1199.5652
$ws.Range("J31").Value = 2193.5715
$ws.Range("L31").Value = 2193.5715
$ws.Range("N31").Value = -2783.5715
# Row 34
$ws.Range("H34").Value = 1199.5652
$ws.Range("J34").Value = 2193.5715
$ws.Range("L34").Value = 2193.5715
$ws.Range("N34").Value = -2597.5715
# Row 132
$ws.Range("H132").Value = 1778.8
$ws.Range("I132").Value = 1322.1904
$ws.Range("J132").Value = 2844.2222
$ws.Range("K132").Value = 3966.5712
$ws.Range("L132").Value = 8532.6666
$ws.Range("M132").Value = -1436.5712
$ws.Range("N132").Value = -13592.6666
# Row 135
$ws.Range("H135").Value = 35503.332
$ws.Range("J135").Value = 35503.332
$ws.Range("L135").Value = 35503.332
$ws.Range("N135").Value = -45643.332

# ================= Sheet: CUL =================
$ws = $wb.Worksheets.Item("CUL")

# Row 2
$ws.Range("H2").Value = 206.33333
$ws.Range("I2").Value = 154.75
$ws.Range("J2").Value = 247.6
$ws.Range("K2").Value = 928.5
$ws.Range("L2").Value = 1485.6
$ws.Range("M2").Value = -815.5
$ws.Range("N2").Value = -1711.6
# Row 5
$ws.Range("H5").Value = 1048.3549
$ws.Range("I5").Value = 1057.1923
$ws.Range("J5").Value = 1002.4
$ws.Range("K5").Value = 3171.5769
$ws.Range("L5").Value = 3007.2
$ws.Range("M5").Value = -3059.5769
$ws.Range("N5").Value = -3231.2
# Row 39
$ws.Range("H39").Value = 4314.857
$ws.Range("J39").Value = 4520.8
$ws.Range("L39").Value = 13562.4
$ws.Range("N39").Value = -14150.4
# Row 55
$ws.Range("H55").Value = 1538.2222
$ws.Range("I55").Value = 67
$ws.Range("J55").Value = 1958.5714
$ws.Range("K55").Value = 201
$ws.Range("L55").Value = 5875.7142
$ws.Range("M55").Value = -24
$ws.Range("N55").Value = -6229.7142
# Row 87
$ws.Range("H87").Value = 847.5
$ws.Range("I87").Value = 847.5
$ws.Range("J87").Value = 0
$ws.Range("K87").Value = 2542.5
$ws.Range("L87").Value = 0
$ws.Range("M87").Value = -1294.5
$ws.Range("N87").ClearContents()
# Row 90
$ws.Range("H90").Value = 847.5
$ws.Range("I90").Value = 847.5
$ws.Range("J90").Value = 0
$ws.Range("K90").Value = 7627.5
$ws.Range("L90").Value = 0
$ws.Range("M90").Value = -1387.5
$ws.Range("N90").ClearContents()
# Row 105
$ws.Range("H105").Value = 117872.11
$ws.Range("J105").Value = 117872.11
$ws.Range("L105").Value = 353616.33
$ws.Range("N105").Value = -358858.33
# Row 107
$ws.Range("H107").Value = 4972.759
$ws.Range("I107").Value = 594
$ws.Range("K107").Value = 1782
$ws.Range("M107").Value = 138
# Row 113
$ws.Range("H113").Value = 680.57574
$ws.Range("J113").Value = 696.0645
$ws.Range("L113").Value = 2088.1935
$ws.Range("N113").Value = -6428.193499999999
# Row 131
$ws.Range("H131").Value = 23291462
$ws.Range("J131").Value = 47732.375
$ws.Range("L131").Value = 143197.125
$ws.Range("N131").Value = -153277.125
# Row 135
$ws.Range("H135").Value = 1048.3549
$ws.Range("I135").Value = 1057.1923
$ws.Range("J135").Value = 1002.4
$ws.Range("K135").Value = 9514.7307
$ws.Range("L135").Value = 9021.6
$ws.Range("M135").Value = -6979.7307
$ws.Range("N135").Value = -14091.6

# ================= Sheet: GSM =================
$ws = $wb.Worksheets.Item("GSM")

# Row 70
$ws.Range("H70").Value = 16670424
$ws.Range("J70").Value = 33336674
$ws.Range("L70").Value = 33336674
$ws.Range("N70").Value = -33337214
# Row 73
$ws.Range("H73").Value = 16670424
$ws.Range("J73").Value = 33336674
$ws.Range("L73").Value = 33336674
$ws.Range("N73").Value = -33338546
# Row 102
$ws.Range("H102").Value = 1408.6945
$ws.Range("I102").Value = 1390.8518
$ws.Range("J102").Value = 1462.2222
$ws.Range("K102").Value = 1390.8518
$ws.Range("L102").Value = 1462.2222
$ws.Range("M102").Value = 231.1482000000001
$ws.Range("N102").Value = -4706.2222
# Row 132
$ws.Range("H132").Value = 2902.2424
$ws.Range("I132").Value = 2537
$ws.Range("J132").Value = 4258.857
$ws.Range("K132").Value = 7611
$ws.Range("L132").Value = 12776.571
$ws.Range("M132").Value = -5081
$ws.Range("N132").Value = -17836.571

# ================= Sheet: LTW =================
$ws = $wb.Worksheets.Item("LTW")

# Row 132
$ws.Range("H132").Value = 3431.3684
$ws.Range("I132").Value = 5749.5
$ws.Range("J132").Value = 2813.2
$ws.Range("K132").Value = 17248.5
$ws.Range("L132").Value = 8439.599999999999
$ws.Range("M132").Value = -14718.5
$ws.Range("N132").Value = -13499.6

# ================= Sheet: WVR =================
$ws = $wb.Worksheets.Item("WVR")

# Row 64
$ws.Range("H64").Value = 17266.666
$ws.Range("J64").Value = 17266.666
$ws.Range("L64").Value = 17266.666
$ws.Range("N64").Value = -17762.666
# Row 67
$ws.Range("H67").Value = 17266.666
$ws.Range("J67").Value = 17266.666
$ws.Range("L67").Value = 17266.666
$ws.Range("N67").Value = -18982.666
# Row 132
$ws.Range("H132").Value = 2917.9656
$ws.Range("I132").Value = 2601.0952
$ws.Range("K132").Value = 7803.285600000001
$ws.Range("M132").Value = -5273.285600000001

